$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to be treated as literal text (these columns store
    # numeric-looking values like "7.00" / "-93.00" / "7.00%" as plain
    # strings, not numbers), then restore the plain "Normal" style so no
    # stray number-format / quote-prefix style is left behind on the cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Insert a new row above the existing row 2 ("Ochieng Charles") so that
# "Jedidah Kemunto" becomes the new row 2 and "Ochieng Charles" shifts to row 3.
# Clear the formatting Excel copies down from the bold header row so the new
# data row stays plain, like the other (unstyled) data rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# New row 2: Jedidah Kemunto
$ws.Range("A2").Value = "Jedidah Kemunto"
Set-TextValue "B2" "0.00"
Set-TextValue "C2" "100.00"
Set-TextValue "D2" "-100.00"
Set-TextValue "E2" "0.00%"

# Insert two new rows after the (now shifted) "Ochieng Charles" row (row 3),
# before the "KD Totals" row, for "Lenah Cheloti" and "Moses  Ngugi".
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).ClearFormats()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).ClearFormats()

# New row 4: Lenah Cheloti
$ws.Range("A4").Value = "Lenah Cheloti"
Set-TextValue "B4" "0.00"
Set-TextValue "C4" "100.00"
Set-TextValue "D4" "-100.00"
Set-TextValue "E4" "0.00%"

# New row 5: Moses  Ngugi
$ws.Range("A5").Value = "Moses  Ngugi"
Set-TextValue "B5" "0.00"
Set-TextValue "C5" "100.00"
Set-TextValue "D5" "-100.00"
Set-TextValue "E5" "0.00%"

# Update the "KD Totals" row, now at row 6, with the new aggregate values.
Set-TextValue "C6" "400.00"
Set-TextValue "D6" "-393.00"
Set-TextValue "E6" "1.75%"
